$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1668.6522
$ws.Range("J112").Value = 1741
$ws.Range("L112").Value = 5223
$ws.Range("N112").Value = -7439
$ws.Range("H124").Value = 48334.75
$ws.Range("J124").Value = 48334.75
$ws.Range("L124").Value = 48334.75
$ws.Range("N124").Value = -58154.75
$ws.Range("H126").Value = 39782.855
$ws.Range("J126").Value = 39782.855
$ws.Range("L126").Value = 39782.855
$ws.Range("N126").Value = -49662.855
$ws.Range("H128").Value = 40661.6
$ws.Range("J128").Value = 40661.6
$ws.Range("L128").Value = 40661.6
$ws.Range("N128").Value = -50621.6
$ws.Range("H129").Value = 1400.6774
$ws.Range("I129").Value = 1607.8889
$ws.Range("J129").Value = 1315.909
$ws.Range("K129").Value = 4823.6667
$ws.Range("L129").Value = 3947.727
$ws.Range("M129").Value = 176.3333000000002
$ws.Range("N129").Value = -13947.727
$ws.Range("H130").Value = 49772
$ws.Range("J130").Value = 49772
$ws.Range("L130").Value = 49772
$ws.Range("N130").Value = -59812
$ws.Range("H133").Value = 35786.92
$ws.Range("J133").Value = 35786.92
$ws.Range("L133").Value = 35786.92
$ws.Range("N133").Value = -45906.92
$ws.Range("H134").Value = 42543.688
$ws.Range("J134").Value = 42543.688
$ws.Range("L134").Value = 42543.688
$ws.Range("N134").Value = -52683.688
$ws.Range("H138").Value = 2476.1914
$ws.Range("I138").Value = 2500.9312
$ws.Range("J138").Value = 2465.1538
$ws.Range("K138").Value = 7502.7936
$ws.Range("L138").Value = 7395.4614
$ws.Range("M138").Value = -2362.7936
$ws.Range("N138").Value = -17675.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9977.203
$ws.Range("I32").Value = 9040.213
$ws.Range("K32").Value = 9040.213
$ws.Range("M32").Value = -8753.213
$ws.Range("H80").Value = 46711.555
$ws.Range("J80").Value = 46711.555
$ws.Range("L80").Value = 46711.555
$ws.Range("N80").Value = -48707.555
$ws.Range("H83").Value = 46711.555
$ws.Range("J83").Value = 46711.555
$ws.Range("L83").Value = 140134.665
$ws.Range("N83").Value = -150118.665
$ws.Range("H123").Value = 35610.5
$ws.Range("J123").Value = 35610.5
$ws.Range("L123").Value = 35610.5
$ws.Range("N123").Value = -45410.5
$ws.Range("H125").Value = 46810
$ws.Range("J125").Value = 46810
$ws.Range("L125").Value = 46810
$ws.Range("N125").Value = -56650
$ws.Range("H130").Value = 42584.2
$ws.Range("J130").Value = 42584.2
$ws.Range("L130").Value = 42584.2
$ws.Range("N130").Value = -52624.2
$ws.Range("H131").Value = 47367.332
$ws.Range("J131").Value = 47367.332
$ws.Range("L131").Value = 47367.332
$ws.Range("N131").Value = -57447.332
$ws.Range("H132").Value = 13890900
$ws.Range("I132").Value = 22728664
$ws.Range("J132").Value = 2984.8572
$ws.Range("K132").Value = 68185992
$ws.Range("L132").Value = 8954.571599999999
$ws.Range("M132").Value = -68183462
$ws.Range("N132").Value = -14014.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 69759.664
$ws.Range("J57").Value = 69759.664
$ws.Range("L57").Value = 69759.664
$ws.Range("N57").Value = -71199.664
$ws.Range("H124").Value = 47997.332
$ws.Range("J124").Value = 47997.332
$ws.Range("L124").Value = 47997.332
$ws.Range("N124").Value = -57817.332
$ws.Range("H125").Value = 50772
$ws.Range("J125").Value = 50772
$ws.Range("L125").Value = 50772
$ws.Range("N125").Value = -60612
$ws.Range("H126").Value = 50768
$ws.Range("J126").Value = 50768
$ws.Range("L126").Value = 50768
$ws.Range("N126").Value = -60648
$ws.Range("H130").Value = 47666
$ws.Range("J130").Value = 47666
$ws.Range("L130").Value = 47666
$ws.Range("N130").Value = -57706
$ws.Range("H132").Value = 44686
$ws.Range("J132").Value = 44686
$ws.Range("L132").Value = 44686
$ws.Range("N132").Value = -54806
$ws.Range("H133").Value = 49499.75
$ws.Range("J133").Value = 49499.75
$ws.Range("L133").Value = 49499.75
$ws.Range("N133").Value = -59619.75
$ws.Range("H134").Value = 2729.6506
$ws.Range("I134").Value = 1421.9062
$ws.Range("J134").Value = 3550.196
$ws.Range("K134").Value = 4265.7186
$ws.Range("L134").Value = 10650.588
$ws.Range("M134").Value = -1730.7186
$ws.Range("N134").Value = -15720.588
$ws.Range("H136").Value = 69759.664
$ws.Range("J136").Value = 69759.664
$ws.Range("L136").Value = 69759.664
$ws.Range("N136").Value = -79959.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H20").Value = 49354.8
$ws.Range("J20").Value = 49354.8
$ws.Range("L20").Value = 49354.8
$ws.Range("N20").Value = -49826.8
$ws.Range("H30").Value = 49354.8
$ws.Range("J30").Value = 49354.8
$ws.Range("L30").Value = 49354.8
$ws.Range("N30").Value = -49536.8
$ws.Range("H31").Value = 2312.09
$ws.Range("I31").Value = 1078.3334
$ws.Range("J31").Value = 3100.8853
$ws.Range("K31").Value = 1078.3334
$ws.Range("L31").Value = 3100.8853
$ws.Range("M31").Value = -783.3334
$ws.Range("N31").Value = -3690.8853
$ws.Range("H34").Value = 2312.09
$ws.Range("I34").Value = 1078.3334
$ws.Range("J34").Value = 3100.8853
$ws.Range("K34").Value = 1078.3334
$ws.Range("L34").Value = 3100.8853
$ws.Range("M34").Value = -876.3334
$ws.Range("N34").Value = -3504.8853
$ws.Range("H58").Value = 1763.5682
$ws.Range("I58").Value = 1528.8572
$ws.Range("J58").Value = 2174.3125
$ws.Range("K58").Value = 1528.8572
$ws.Range("L58").Value = 2174.3125
$ws.Range("M58").Value = -1325.8572
$ws.Range("N58").Value = -2580.3125
$ws.Range("H128").Value = 49354.8
$ws.Range("J128").Value = 49354.8
$ws.Range("L128").Value = 49354.8
$ws.Range("N128").Value = -59314.8
$ws.Range("H132").Value = 39729.055
$ws.Range("I132").Value = 1630.3448
$ws.Range("J132").Value = 177836.88
$ws.Range("K132").Value = 4891.0344
$ws.Range("L132").Value = 533510.64
$ws.Range("M132").Value = -2361.0344
$ws.Range("N132").Value = -538570.64
$ws.Range("H135").Value = 45080
$ws.Range("J135").Value = 45080
$ws.Range("L135").Value = 45080
$ws.Range("N135").Value = -55220
$ws.Range("H136").Value = 1763.5682
$ws.Range("I136").Value = 1528.8572
$ws.Range("J136").Value = 2174.3125
$ws.Range("K136").Value = 4586.571599999999
$ws.Range("L136").Value = 6522.9375
$ws.Range("M136").Value = -2036.571599999999
$ws.Range("N136").Value = -11622.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2443.9832
$ws.Range("I122").Value = 699.08105
$ws.Range("J122").Value = 5378.591
$ws.Range("K122").Value = 6291.72945
$ws.Range("L122").Value = 48407.319
$ws.Range("M122").Value = -3841.72945
$ws.Range("N122").Value = -53307.319
$ws.Range("H133").Value = 5896.8696
$ws.Range("J133").Value = 5679.9
$ws.Range("L133").Value = 17039.7
$ws.Range("N133").Value = -27159.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1434.1578
$ws.Range("I113").Value = 1657.1428
$ws.Range("J113").Value = 1304.0834
$ws.Range("K113").Value = 1657.1428
$ws.Range("L113").Value = 1304.0834
$ws.Range("M113").Value = 512.8571999999999
$ws.Range("N113").Value = -5644.0834
$ws.Range("H130").Value = 44956
$ws.Range("J130").Value = 44956
$ws.Range("L130").Value = 44956
$ws.Range("N130").Value = -54996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 48707
$ws.Range("J36").Value = 48707
$ws.Range("L36").Value = 48707
$ws.Range("N36").Value = -49831
$ws.Range("H81").Value = 42181
$ws.Range("J81").Value = 42181
$ws.Range("L81").Value = 42181
$ws.Range("N81").Value = -44177
$ws.Range("H84").Value = 42181
$ws.Range("J84").Value = 42181
$ws.Range("L84").Value = 126543
$ws.Range("N84").Value = -136527
$ws.Range("H130").Value = 39476.332
$ws.Range("J130").Value = 39476.332
$ws.Range("L130").Value = 39476.332
$ws.Range("N130").Value = -49516.332
$ws.Range("H136").Value = 1657.9688
$ws.Range("I136").Value = 1380.091
$ws.Range("J136").Value = 2269.3
$ws.Range("K136").Value = 4140.272999999999
$ws.Range("L136").Value = 6807.900000000001
$ws.Range("M136").Value = -1590.272999999999
$ws.Range("N136").Value = -11907.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 47420
$ws.Range("J120").Value = 47420
$ws.Range("L120").Value = 47420
$ws.Range("N120").Value = -57096
$ws.Range("H128").Value = 48895
$ws.Range("J128").Value = 48895
$ws.Range("L128").Value = 48895
$ws.Range("N128").Value = -58855
$ws.Range("H138").Value = 45387.5
$ws.Range("J138").Value = 45387.5
$ws.Range("L138").Value = 45387.5
$ws.Range("N138").Value = -55667.5
